$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set cell C6 to "YES" (was "NO")
$ws.Range("C6").Value = "YES"

# Update selection to just C6 (was C2:C9 with active cell C2)
$ws.Range("C6").Select()
